# Apply updated dSF (column F) values to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F3").Value = -9
$ws.Range("F4").Value = -10
$ws.Range("F5").Value = -1
$ws.Range("F6").Value = -3
$ws.Range("F10").Value = -4
$ws.Range("F14").Value = -4
$ws.Range("F17").Value = -1
$ws.Range("F21").Value = -1
$ws.Range("F22").Value = -1
